# Update "想去人数" (interest count) figures for several events that
# appear both on the "展览" sheet and on the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibit.Range("F3").Value  = 4155
$wsExhibit.Range("F8").Value  = 42
$wsExhibit.Range("F9").Value  = 208
$wsExhibit.Range("F11").Value = 110
$wsExhibit.Range("F13").Value = 1556
$wsExhibit.Range("F14").Value = 285
$wsExhibit.Range("F15").Value = 3113

# --- 全部类型 (sheet4) ---
$wsAll.Range("F3").Value  = 4155
$wsAll.Range("F9").Value  = 42
$wsAll.Range("F11").Value = 208
$wsAll.Range("F13").Value = 110
$wsAll.Range("F17").Value = 1556
$wsAll.Range("F18").Value = 285
$wsAll.Range("F19").Value = 3113
